$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    11 = 17.517
    12 = 17.386
    15 = 16.357
    27 = 16.899
    28 = 17.148
    31 = 16.509
    32 = 17.517
    36 = 16.652
    38 = 16.631
    46 = 17.015
    54 = 16.754
    55 = 16.378
    56 = 16.351
    67 = 17.24
    69 = 17.386
    72 = 16.898
    73 = 16.602
    83 = 16.822
    86 = 16.542
    91 = 17.577
    93 = 17.22
    99 = 16.766
}

foreach ($row in $updates.Keys) {
    $ws.Range("E$row").Value = $updates[$row]
}
